# edit.ps1 - Apply the HealthCardProject.docx revision:
#   1. Title "Medical Health Card-main" -> "Medical Health Card", color FF0000 -> 4472C4/accent1
#   2. Add "Sample file" run to the first empty (sz=32) paragraph after "Overview/Introduction:"
#   3. Merge the split "Database that should hold the history of users treatment" runs
#      (dropping the gramStart/gramEnd proofErr markers) into a single run
#   4. Merge the split "Various ways to put data into database. These may be :" runs
#      (dropping the gramStart/gramEnd proofErr markers) into a single run

$d = $word.ActiveDocument

$pkgNs = 'xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"'
$wNs   = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function New-WordPackageXml([string]$bodyXml) {
    return '<pkg:package ' + $pkgNs + '><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1. Title paragraph: recolor to theme accent1 and drop the "-main" run
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range

$titleBody = '<w:p w14:paraId="4FD8F73A" w14:textId="798996AF" w:rsidR="003478AA" w:rsidRPr="00C4471D" w:rsidRDefault="00D977AA" w:rsidP="003478AA">' +
    '<w:pPr><w:pStyle w:val="Title"/><w:ind w:left="2160" w:firstLine="720"/>' +
    '<w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00C4471D"><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:u w:val="single"/></w:rPr>' +
    '<w:t xml:space="preserve">Medical </w:t></w:r>' +
    '<w:r w:rsidR="003478AA" w:rsidRPr="00C4471D"><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:u w:val="single"/></w:rPr>' +
    '<w:t>Health Card</w:t></w:r>' +
    '</w:p>'

$titleRange.InsertXML((New-WordPackageXml $titleBody))

# ---------------------------------------------------------------------------
# 2. First empty paragraph after the heading gets a "Sample file" run
# ---------------------------------------------------------------------------
$sampleRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq [string][char]13 -and $i -gt 3) {
        $sampleRange = $cand.Range
        break
    }
}
if ($sampleRange -eq $null) {
    $sampleRange = $d.Paragraphs(4).Range
}

$sampleBody = '<w:p w14:paraId="3A3532FF" w14:textId="25F23D1A" w:rsidR="003478AA" w:rsidRDefault="003478AA">' +
    '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Sample file</w:t></w:r>' +
    '</w:p>'

$sampleRange.InsertXML((New-WordPackageXml $sampleBody))

# ---------------------------------------------------------------------------
# 3 & 4. Merge the grammar-split runs back into single runs, dropping proofErr
# ---------------------------------------------------------------------------
function Find-ParagraphIndexByText([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

$dbIdx = Find-ParagraphIndexByText "Database that should hold"
$dbRange = $d.Paragraphs($dbIdx).Range
$dbBody = '<w:p w14:paraId="24BB652A" w14:textId="77777777" w:rsidR="003478AA" w:rsidRPr="003478AA" w:rsidRDefault="003478AA" w:rsidP="003478AA">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr>' +
    '<w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="003478AA"><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t>Database that should hold the history of users treatment</w:t></w:r>' +
    '</w:p>'
$dbRange.InsertXML((New-WordPackageXml $dbBody))

$wayIdx = Find-ParagraphIndexByText "Various ways to put data"
$wayRange = $d.Paragraphs($wayIdx).Range
$wayBody = '<w:p w14:paraId="4FF1FFE6" w14:textId="77777777" w:rsidR="003478AA" w:rsidRPr="003478AA" w:rsidRDefault="003478AA" w:rsidP="003478AA">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr>' +
    '<w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="003478AA"><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t>Various ways to put data into database. These may be :</w:t></w:r>' +
    '</w:p>'
$wayRange.InsertXML((New-WordPackageXml $wayBody))

Write-Output "Edits applied successfully"
